$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptocurrency price / volume snapshot values.
# Column D ("Price") cells are forced to Text format ("@") before the value
# is written so that Excel keeps the exact original text representation
# (thousand-dot grouping, trailing zeros, etc.) instead of re-parsing the
# string into a floating point number.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '46.670.48'
$ws.Range("E2").Value = '  -0.08%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.276.62'
$ws.Range("E3").Value = '  -2.35%  '

$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '301.45'
$ws.Range("E5").Value = '  -1.42%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '100.28'
$ws.Range("E6").Value = '  +2.37%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.566'
$ws.Range("E7").Value = '  -1.71%  '

$ws.Range("E8").Value = '  -0.07%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.510'
$ws.Range("E9").Value = '  -5.38%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.30'
$ws.Range("E10").Value = '  -2.14%  '

$ws.Range("E11").Value = '  -0.88%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.12'
$ws.Range("E12").Value = '  -4.52%  '

$ws.Range("E13").Value = '  -1.60%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.620.84'

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.271.80'
$ws.Range("E15").Value = '  -2.82%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '13.70'
$ws.Range("E16").Value = '  -3.21%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.804'
$ws.Range("E17").Value = '  -3.51%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '46.622.17'
$ws.Range("E18").Value = '  +0.03%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.65'
$ws.Range("E19").Value = '  -3.13%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0971'
$ws.Range("E20").Value = '  +2.57%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.88'
$ws.Range("E21").Value = '  -5.03%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '66.04'
$ws.Range("E22").Value = '  -1.22%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '249.60'
$ws.Range("E23").Value = '  +1.85%  '

$ws.Range("E24").Value = '  -5.37%  '

$ws.Range("E25").Value = '  +0.11%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.88'
$ws.Range("E26").Value = '  -5.13%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '41.66'
$ws.Range("E27").Value = '  -1.18%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.22'
$ws.Range("E28").Value = '  -2.92%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.67'
$ws.Range("E29").Value = '  -1.81%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '20.24'
$ws.Range("E30").Value = '  +0.24%  '

$ws.Range("E31").Value = '  +7.33%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.42'
$ws.Range("E32").Value = '  +13.54%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '147.31'
$ws.Range("E33").Value = '  -3.30%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.39'
$ws.Range("E34").Value = '  -5.98%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0774'
$ws.Range("E35").Value = '  -4.18%  '

$ws.Range("E36").Value = '  +7.05%  '

$ws.Range("E37").Value = '  -2.81%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '15.89'
$ws.Range("E38").Value = '  +13.49%  '

$ws.Range("E39").Value = '  -5.80%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.91'
$ws.Range("E40").Value = '  -2.75%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0298'
$ws.Range("E41").Value = '  -6.33%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.15'
$ws.Range("E42").Value = '  -7.59%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.999'
$ws.Range("E43").Value = '  -0.04%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '93.27'
$ws.Range("E44").Value = '  +15.52%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.798.18'
$ws.Range("E45").Value = '  +0.51%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.89'
$ws.Range("E46").Value = '  -2.46%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '71.46'
$ws.Range("E47").Value = '  -3.98%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.187'
$ws.Range("E48").Value = '  -6.25%  '

$ws.Range("E49").Value = '  -1.46%  '

$ws.Range("B50").Value = 'FraxShare'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.94'
$ws.Range("E50").Value = '  -0.97%  '

$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '95.25'
$ws.Range("E51").Value = '  -3.04%  '
